$wb = $excel.ActiveWorkbook

# "想去人数" (F column) head-count estimates were refreshed for this gh-pages data rebuild.
# Both the "展览" (Exhibitions) sheet and the "全部类型" (All Types) rollup sheet
# carry the same 24 rows of updated figures.
$sheetNames = @("展览", "全部类型")

$updates = @{
    2 = 1038
    3 = 720
    4 = 252
    6 = 1090
    8 = 1639
    9 = 6032
    10 = 472
    11 = 334
    12 = 267
    13 = 80
    14 = 355
    16 = 4756
    17 = 254
    18 = 1256
    21 = 216
    22 = 94
    23 = 241
    24 = 90
    26 = 2
    31 = 75
    32 = 36
    33 = 52
    34 = 18
    35 = 56
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

Write-Output "updated F column on $($sheetNames.Count) sheets"
